# PartList.xlsx update: add "Mechanical Parts" + "Small Parts" sections
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header: A1 gets the "Good" (Gut) cell style --------------------
$ws.Range("A1").Style = "Gut"

# --- New section starting at row 25 -----------------------------------------
$ws.Range("A25").Value = "Mechanical Parts (Self Print etc)"
$ws.Range("A25").Style = "Gut"

# Table header row 26 (bold header style, same as row 3)
$ws.Range("A26").Value = "Name"
$ws.Range("B26").Value = "File Location"
$ws.Range("C26").Value = "How to manifacture"
$ws.Range("D26").Value = "Amount"
$ws.Range("A26:D26").Font.Bold = $true

# I26 uses the "Good" style too (empty cell, formatting only)
$ws.Range("I26").Style = "Gut"

# Row 27: Base_Protection
$ws.Range("A27").Value = "Base_Protection"
$ws.Range("B27").Value = "112_CAD Files\Base"
$ws.Range("C27").Value = "* 3D Print or `n* Laser Cutting or `n* manual (sawing/drilling)"
$ws.Range("D27").Value = 2
$ws.Range("A27:D27").VerticalAlignment = -4160
$ws.Range("C27").WrapText = $true
$ws.Rows(27).RowHeight = 45

# Row 28: Hallsensor_DistanceBolt_8mm4
$ws.Range("A28").Value = "Hallsensor_DistanceBolt_8mm4"
$ws.Range("B28").Value = "112_CAD Files\SensorHead"
$ws.Range("C28").Value = "* 3D Print (PETG)"
$ws.Range("D28").Value = 3
$ws.Range("A28:D28").VerticalAlignment = -4160

# Row 29 stays blank but keeps the vertical-top formatting from the selection
$ws.Range("A29:B29").VerticalAlignment = -4160

# Row 30: Small Parts sub-header (Good style + vertical top)
$ws.Range("A30").Value = "Small Parts"
$ws.Range("A30").Style = "Gut"
$ws.Range("A30").VerticalAlignment = -4160

# Row 31: second mini table header
$ws.Range("A31").Value = "Name"
$ws.Range("B31").Value = "Amount"
$ws.Range("A31:B31").Font.Bold = $true
$ws.Range("A31:B31").VerticalAlignment = -4160

# Row 32
$ws.Range("A32").Value = "M3 15mm Plastic Screws with washer/nut"
$ws.Range("B32").Value = 3
$ws.Range("A32:B32").VerticalAlignment = -4160

# Row 33
$ws.Range("A33").Value = "Cable tie 2.5mm"
$ws.Range("B33").Value = 2
$ws.Range("A33:B33").VerticalAlignment = -4160

# Rows 34-45 stay blank, still carrying the vertical-top formatting
$ws.Range("A34:B45").VerticalAlignment = -4160

# --- Selection cursor, matching the saved UI state --------------------------
$ws.Range("B38").Select() | Out-Null
